# Auto-generated edit script for commit 'Add data for 2022-08-25'
# Applies 85 cell updates (mostly +1 increments reflecting newly-ingested
# crime records for 2022-08-25) across 19 worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("B2").Value = 28
$ws.Range("D2").Value = 64
$ws.Range("E3").Value = 96
$ws.Range("F3").Value = 92
$ws.Range("I3").Value = 133
$ws.Range("B6").Value = 255
$ws.Range("C6").Value = 325
$ws.Range("D6").Value = 294
$ws.Range("E6").Value = 286
$ws.Range("F6").Value = 378
$ws.Range("G6").Value = 332
$ws.Range("B7").Value = 348
$ws.Range("C7").Value = 437
$ws.Range("D7").Value = 457
$ws.Range("E7").Value = 436
$ws.Range("F7").Value = 537
$ws.Range("G7").Value = 482
$ws.Range("I7").Value = 592

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("D6").Value = 30
$ws.Range("D7").Value = 39

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("B6").Value = 13
$ws.Range("D6").Value = 16
$ws.Range("F6").Value = 16
$ws.Range("G6").Value = 8
$ws.Range("B7").Value = 18
$ws.Range("D7").Value = 26
$ws.Range("F7").Value = 36
$ws.Range("G7").Value = 16

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("C4").Value = 5
$ws.Range("C5").Value = 7

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("C6").Value = 33
$ws.Range("C7").Value = 37

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("F7").Value = 8
$ws.Range("E8").Value = 32
$ws.Range("G19").Value = 18
$ws.Range("B21").Value = 5
$ws.Range("C28").Value = 37
$ws.Range("D32").Value = 39
$ws.Range("B36").Value = 18
$ws.Range("D36").Value = 26
$ws.Range("F36").Value = 36
$ws.Range("G36").Value = 16
$ws.Range("F47").Value = 12
$ws.Range("I53").Value = 93
$ws.Range("D62").Value = 2
$ws.Range("D76").Value = 10
$ws.Range("E77").Value = 21
$ws.Range("B87").Value = 4
$ws.Range("C88").Value = 7
$ws.Range("F92").Value = 6
$ws.Range("B94").Value = 3
$ws.Range("B96").Value = 12
$ws.Range("B98").Value = 348
$ws.Range("C98").Value = 437
$ws.Range("D98").Value = 457
$ws.Range("E98").Value = 436
$ws.Range("F98").Value = 537
$ws.Range("G98").Value = 482
$ws.Range("I98").Value = 592

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("B2").Value = 2
$ws.Range("B6").Value = 12

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I3").Value = 21
$ws.Range("I7").Value = 93

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 6

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("D2").Value = 2
$ws.Range("D6").Value = 10

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = 4

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("D4").Value = 2
$ws.Range("D5").Value = 2

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("E3").Value = 5
$ws.Range("E7").Value = 21

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("G5").Value = 12
$ws.Range("G6").Value = 18

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("F3").Value = 2
$ws.Range("F6").Value = 8

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("F3").Value = 2
$ws.Range("F6").Value = 12

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 3

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("E5").Value = 24
$ws.Range("E6").Value = 32

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 5
